# "Add hashtable use to parameters" - add a runmode (Y/N) column driving
# which data rows / test cases are exercised, plus a couple of data fixes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: AddCustomerTest
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("AddCustomerTest")

# fix header typo
$ws1.Range("A1").Value = "firstname"

# new "runmode" header + column
$ws1.Range("E1").Value = "runmode"

# existing rows get a runmode flag
$ws1.Range("E2").Value = "Y"
$ws1.Range("E3").Value = "Y"
$ws1.Range("E4").Value = "N"

# row 4 name correction: Enzo -> Jose
$ws1.Range("A4").Value = "Jose"

# brand new data row
$ws1.Range("A5").Value = "Jorge"
$ws1.Range("B5").Value = "Souza"
$ws1.Range("C5").Value = 789456
$ws1.Range("D5").Value = "Customer added successfully"
$ws1.Range("E5").Value = "N"

$ws1.Columns.Item(5).ColumnWidth = 10.8

[void]$ws1.Range("E1").Select()

# ---------------------------------------------------------------------
# Sheet 3: test_suite - run every test now (was N for OpenAccountTest)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("test_suite")

$ws3.Range("B2").Value = "Y"
$ws3.Range("B3").Value = "Y"
$ws3.Range("B4").Value = "Y"

[void]$ws3.Range("B3").Select()

# ---------------------------------------------------------------------
# Sheet 2: OpenAccountTest (ends up the active tab)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("OpenAccountTest")

$ws2.Range("C1").Value = "runmode"
$ws2.Range("C2").Value = "Y"

$ws2.Columns.Item(3).ColumnWidth = 10.8

$ws2.Activate()
[void]$ws2.Range("C2").Select()
